# Update '想去人数' (interested count) figures in column F across sheets
# as published in the refreshed 苏州-漫展信息 data snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 11950  # was 11945
$ws.Range("F4").Value = 29  # was 26
$ws.Range("F5").Value = 225  # was 224
$ws.Range("F8").Value = 11846  # was 11839
$ws.Range("F9").Value = 496  # was 494
$ws.Range("F12").Value = 70  # was 69
$ws.Range("F13").Value = 1786  # was 1785
$ws.Range("F14").Value = 5869  # was 5867
$ws.Range("F15").Value = 128  # was 126
$ws.Range("F16").Value = 3543  # was 3542
$ws.Range("F17").Value = 191  # was 190
$ws.Range("F18").Value = 25  # was 23

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 574  # was 575

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 574  # was 575
$ws.Range("F5").Value = 11950  # was 11945
$ws.Range("F6").Value = 29  # was 26
$ws.Range("F7").Value = 225  # was 224
$ws.Range("F11").Value = 11846  # was 11839
$ws.Range("F12").Value = 496  # was 494
$ws.Range("F15").Value = 70  # was 69
$ws.Range("F16").Value = 1786  # was 1785
$ws.Range("F18").Value = 5869  # was 5868
$ws.Range("F19").Value = 128  # was 126
$ws.Range("F20").Value = 3543  # was 3542
$ws.Range("F21").Value = 191  # was 190
$ws.Range("F22").Value = 25  # was 23
